$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 4 and 5 (pushes the existing
# rows 4-29 down to rows 6-31, which already carries all of their data
# and labels along for free).
$ws.Rows("4:5").Insert()

# The insert leaves A4/A5 with a stray auto-generated style; restore the
# same bordered/bold/centered style used by every other cell in column A.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the two brand-new simulation rows: "Holden" and "Rizzie Spiral".
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$row4 = @(0.982283449949899, 1.004341905799651, 0.982283449949899, 1.004477599685333, 1.004438038641642, 0.9902244938317223, 1.011966116881022, 1.004341905799651, 1.004341905799651, 1.004477599685333, 0.9933805248176162, 0.9933805248176162, 0.9923285144889848, 0.9970343184782943, 0.9970343184782943, 0.9988612153086334, 0.9988612153086334, 0.9996219341315449)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4[$i]
}

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$row5 = @(0.9877769389850531, 1.003010476454693, 0.9877769389850531, 1.003080927444437, 1.003060383958551, 0.9932610206827857, 1.008229223149272, 1.003010476454693, 1.003010476454693, 1.003080927444437, 0.995428933214745, 0.995428933214745, 0.9947062957040919, 0.9979561142947277, 0.9979561142947279, 0.9992197048347193, 0.9992197048347193, 0.9997364951124653)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5[$i]
}

# Rename "Thomas Hex" (now sitting two rows further down) to "Matthies Hex".
$ws.Range("B11").Value = "Matthies Hex"

Write-Host "edit complete"
